# Update Leve profit figures across sheets (scheduled runner refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 78123.766
$ws.Range("J93").Value = 78123.766
$ws.Range("L93").Value = 78123.766
$ws.Range("N93").Value = -83115.766

$ws.Range("H101").Value = 3266
$ws.Range("I101").Value = 1603.4286
$ws.Range("K101").Value = 4810.2858
$ws.Range("M101").Value = -3188.2858

$ws.Range("H124").Value = 37245.383
$ws.Range("J124").Value = 37245.383
$ws.Range("L124").Value = 37245.383
$ws.Range("N124").Value = -47065.383

$ws.Range("H130").Value = 45406.668
$ws.Range("J130").Value = 45406.668
$ws.Range("L130").Value = 45406.668
$ws.Range("N130").Value = -55446.668

$ws.Range("H132").Value = 2633.7021
$ws.Range("I132").Value = 1276.6571
$ws.Range("J132").Value = 6591.75
$ws.Range("K132").Value = 3829.9713
$ws.Range("L132").Value = 19775.25
$ws.Range("M132").Value = -1299.9713
$ws.Range("N132").Value = -24835.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 25499.6
$ws.Range("J44").Value = 25499.6
$ws.Range("L44").Value = 25499.6
$ws.Range("N44").Value = -26475.6

$ws.Range("H55").Value = 24640.834
$ws.Range("J55").Value = 24640.834
$ws.Range("L55").Value = 24640.834
$ws.Range("N55").Value = -25270.834

$ws.Range("H80").Value = 25652.092
$ws.Range("J80").Value = 25652.092
$ws.Range("L80").Value = 25652.092
$ws.Range("N80").Value = -27648.092

$ws.Range("H83").Value = 25652.092
$ws.Range("J83").Value = 25652.092
$ws.Range("L83").Value = 76956.276
$ws.Range("N83").Value = -86940.276

$ws.Range("H129").Value = 49099.6
$ws.Range("J129").Value = 49099.6
$ws.Range("L129").Value = 49099.6
$ws.Range("N129").Value = -59099.6

$ws.Range("H133").Value = 68000
$ws.Range("J133").Value = 68000
$ws.Range("L133").Value = 68000
$ws.Range("N133").Value = -73060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 9367.5
$ws.Range("J51").Value = 10990
$ws.Range("L51").Value = 10990
$ws.Range("N51").Value = -12462

$ws.Range("H61").Value = 9367.5
$ws.Range("J61").Value = 10990
$ws.Range("L61").Value = 10990
$ws.Range("N61").Value = -11686

$ws.Range("H123").Value = 52653.332
$ws.Range("J123").Value = 52653.332
$ws.Range("L123").Value = 52653.332
$ws.Range("N123").Value = -62453.332

$ws.Range("H134").Value = 4635.242
$ws.Range("I134").Value = 5106.0356
$ws.Range("K134").Value = 15318.1068
$ws.Range("M134").Value = -12783.1068

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1170.5714
$ws.Range("I5").Value = 848.5
$ws.Range("K5").Value = 2545.5
$ws.Range("M5").Value = -2433.5

$ws.Range("H129").Value = 57024.89
$ws.Range("I129").Value = 91715.45
$ws.Range("J129").Value = 2511.1428
$ws.Range("K129").Value = 275146.35
$ws.Range("L129").Value = 7533.428400000001
$ws.Range("M129").Value = -270146.35
$ws.Range("N129").Value = -17533.4284

$ws.Range("H130").Value = 2000
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").Value = $null

$ws.Range("H131").Value = 846.1667
$ws.Range("I131").Value = 357.27274
$ws.Range("J131").Value = 934.3279
$ws.Range("K131").Value = 1071.81822
$ws.Range("L131").Value = 2802.9837
$ws.Range("M131").Value = 3968.18178
$ws.Range("N131").Value = -12882.9837

$ws.Range("H132").Value = 1648009.8
$ws.Range("I132").Value = 2632815.5
$ws.Range("J132").Value = 6666.6665
$ws.Range("K132").Value = 23695339.5
$ws.Range("L132").Value = 59999.9985
$ws.Range("M132").Value = -23692809.5
$ws.Range("N132").Value = -65059.9985

$ws.Range("H135").Value = 1170.5714
$ws.Range("I135").Value = 848.5
$ws.Range("K135").Value = 7636.5
$ws.Range("M135").Value = -5101.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 19076.666
$ws.Range("J62").Value = 19076.666
$ws.Range("L62").Value = 19076.666
$ws.Range("N62").Value = -20448.666

$ws.Range("H65").Value = 19076.666
$ws.Range("J65").Value = 19076.666
$ws.Range("L65").Value = 57229.99800000001
$ws.Range("N65").Value = -64093.99800000001

$ws.Range("H93").Value = 9692.200000000001
$ws.Range("J93").Value = 9692.200000000001
$ws.Range("L93").Value = 9692.200000000001
$ws.Range("N93").Value = -13436.2

$ws.Range("H102").Value = 1431.8
$ws.Range("I102").Value = 1285.2778
$ws.Range("J102").Value = 1808.5714
$ws.Range("K102").Value = 1285.2778
$ws.Range("L102").Value = 1808.5714
$ws.Range("M102").Value = 336.7221999999999
$ws.Range("N102").Value = -5052.5714

$ws.Range("H124").Value = 49716
$ws.Range("J124").Value = 49716
$ws.Range("L124").Value = 49716
$ws.Range("N124").Value = -59536

$ws.Range("H127").Value = 39985
$ws.Range("J127").Value = 39985
$ws.Range("L127").Value = 39985
$ws.Range("N127").Value = -49905

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 20625.5
$ws.Range("J108").Value = 20625.5
$ws.Range("L108").Value = 20625.5
$ws.Range("N108").Value = -28305.5

$ws.Range("H123").Value = 38599.832
$ws.Range("J123").Value = 38599.832
$ws.Range("L123").Value = 38599.832
$ws.Range("N123").Value = -48399.832

$ws.Range("H128").Value = 50323.168
$ws.Range("J128").Value = 50323.168
$ws.Range("L128").Value = 50323.168
$ws.Range("N128").Value = -60283.168

$ws.Range("H130").Value = 37335.8
$ws.Range("J130").Value = 37335.8
$ws.Range("L130").Value = 37335.8
$ws.Range("N130").Value = -47375.8

$ws.Range("H132").Value = 2946.3845
$ws.Range("I132").Value = 2142.8572
$ws.Range("J132").Value = 3883.8333
$ws.Range("K132").Value = 6428.571599999999
$ws.Range("L132").Value = 11651.4999
$ws.Range("M132").Value = -3898.571599999999
$ws.Range("N132").Value = -16711.4999

$ws.Range("H134").Value = 28746
$ws.Range("J134").Value = 28746
$ws.Range("L134").Value = 28746
$ws.Range("N134").Value = -38886

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2000
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = $null

$ws.Range("H109").Value = 25888.5
$ws.Range("J109").Value = 25888.5
$ws.Range("L109").Value = 25888.5
$ws.Range("N109").Value = -28662.5

$ws.Range("H129").Value = 25691.215
$ws.Range("J129").Value = 25691.215
$ws.Range("L129").Value = 25691.215
$ws.Range("N129").Value = -35691.215

$ws.Range("H132").Value = 3347
$ws.Range("I132").Value = 3928.5
$ws.Range("J132").Value = 2818.3635
$ws.Range("K132").Value = 11785.5
$ws.Range("L132").Value = 8455.0905
$ws.Range("M132").Value = -9255.5
$ws.Range("N132").Value = -13515.0905
